$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the command token text to match new mixed-case token reading
$ws.Range("A1").Value = "PARAgrAPH"

# Move selection to C5 as last done by the user
$ws.Range("C5").Select()
